$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AA1").Value = "צוללת"
$ws.Range("AA8").Value = "עניבה"
$ws.Range("AA11").Value = "קלפים"
$ws.Range("AA14").Value = "כביסה"
$ws.Range("AA16").Value = "משאבה"
$ws.Range("AA20").Value = "מייבש"
$ws.Range("AA22").Value = "מנורה"
$ws.Range("AA23").Value = "בלנדר"
$ws.Range("AA27").Value = "גיטרה"
$ws.Range("A28").Value = "צוללת"
$ws.Range("AB28").Value = "צוללת"
$ws.Range("AC28").Value = "צוללת"
$ws.Range("AD28").Value = "צוללת"
$ws.Range("B28").Value = "צוללת"
$ws.Range("C28").Value = "צוללת"
$ws.Range("F28").Value = "צוללת"
$ws.Range("G28").Value = "צוללת"
$ws.Range("I28").Value = "צוללת"
$ws.Range("J28").Value = "צוללת"
$ws.Range("K28").Value = "צוללת"
$ws.Range("M28").Value = "צוללת"
$ws.Range("O28").Value = "צוללת"
$ws.Range("Q28").Value = "צוללת"
$ws.Range("R28").Value = "צוללת"
$ws.Range("S28").Value = "צוללת"
$ws.Range("T28").Value = "צוללת"
$ws.Range("U28").Value = "צוללת"
$ws.Range("V28").Value = "צוללת"
$ws.Range("Z28").Value = "צוללת"
$ws.Range("AA29").Value = "מגירה"
$ws.Range("AA30").Value = "יאכטה"
$ws.Range("AA31").Value = "ספינה"

$ws.Range("AA5").ClearContents()
$ws.Range("AA6").ClearContents()
$ws.Range("AA13").ClearContents()
$ws.Range("AA15").ClearContents()
$ws.Range("AA17").ClearContents()
$ws.Range("AA24").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("L28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("P28").ClearContents()
$ws.Range("W28").ClearContents()
